# Insert a new weekly price record for "Ajo" (Vega Monumental Concepción)
# right before the existing row 159, shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 159..196 down to 160..197 and open up a blank row 159.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new weekly observation.
$ws.Cells.Item(159, 1).Value  = 11
$ws.Cells.Item(159, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(159, 3).Value  = "Bíobío"
$ws.Cells.Item(159, 4).Value  = 44782
$ws.Cells.Item(159, 5).Value  = 8
$ws.Cells.Item(159, 6).Value  = 100112003
$ws.Cells.Item(159, 7).Value  = "Ajo"
$ws.Cells.Item(159, 8).Value  = "Chino"
$ws.Cells.Item(159, 9).Value  = "Primera"
$ws.Cells.Item(159, 10).Value = 400
$ws.Cells.Item(159, 11).Value = 23000
$ws.Cells.Item(159, 12).Value = 24000
$ws.Cells.Item(159, 13).Value = 23500
$ws.Cells.Item(159, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(159, 15).Value = "China"
$ws.Cells.Item(159, 16).Value = 2350
$ws.Cells.Item(159, 17).Value = 10
$ws.Cells.Item(159, 18).Value = "Hortaliza"
